# Update "想去人数" (wanted-to-go count) values in column F across all four
# worksheets to reflect the regenerated gh-pages data snapshot.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1617
$ws.Range("F6").Value = 794
$ws.Range("F7").Value = 696
$ws.Range("F8").Value = 1287
$ws.Range("F9").Value = 2609
$ws.Range("F10").Value = 1345
$ws.Range("F11").Value = 476
$ws.Range("F12").Value = 2315
$ws.Range("F13").Value = 2038
$ws.Range("F14").Value = 713
$ws.Range("F15").Value = 6450
$ws.Range("F17").Value = 1227
$ws.Range("F18").Value = 138
$ws.Range("F19").Value = 1466
$ws.Range("F20").Value = 1330
$ws.Range("F21").Value = 1192
$ws.Range("F23").Value = 2268
$ws.Range("F24").Value = 1107
$ws.Range("F25").Value = 728
$ws.Range("F26").Value = 237
$ws.Range("F27").Value = 5299
$ws.Range("F28").Value = 284
$ws.Range("F29").Value = 1252
$ws.Range("F30").Value = 45
$ws.Range("F31").Value = 3705
$ws.Range("F32").Value = 636
$ws.Range("F33").Value = 1685
$ws.Range("F34").Value = 1074
$ws.Range("F35").Value = 158
$ws.Range("F36").Value = 273
$ws.Range("F37").Value = 960
$ws.Range("F38").Value = 1047
$ws.Range("F39").Value = 390
$ws.Range("F40").Value = 1766
$ws.Range("F42").Value = 104
$ws.Range("F43").Value = 896
$ws.Range("F44").Value = 1046
$ws.Range("F45").Value = 510
$ws.Range("F49").Value = 77

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 440
$ws.Range("F9").Value = 491
$ws.Range("F11").Value = 386
$ws.Range("F15").Value = 964
$ws.Range("F18").Value = 14
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = 78
$ws.Range("F21").Value = 589
$ws.Range("F22").Value = 237
$ws.Range("F23").Value = 355
$ws.Range("F25").Value = 171
$ws.Range("F31").Value = 35
$ws.Range("F35").Value = 44
$ws.Range("F38").Value = 193
$ws.Range("F43").Value = 7

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3295
$ws.Range("F7").Value = 1456
$ws.Range("F8").Value = 764
$ws.Range("F9").Value = 386
$ws.Range("F10").Value = 2779
$ws.Range("F11").Value = 285
$ws.Range("F12").Value = 524
$ws.Range("F13").Value = 554
$ws.Range("F14").Value = 1152

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1456
$ws.Range("F3").Value = 764
$ws.Range("F5").Value = 386
$ws.Range("F6").Value = 2779
$ws.Range("F7").Value = 1617
$ws.Range("F8").Value = 794
$ws.Range("F9").Value = 696
$ws.Range("F10").Value = 1287
$ws.Range("F11").Value = 2609
$ws.Range("F12").Value = 1345
$ws.Range("F13").Value = 476
$ws.Range("F14").Value = 2315
$ws.Range("F15").Value = 2038
$ws.Range("F16").Value = 713
$ws.Range("F17").Value = 6450
$ws.Range("F19").Value = 524
$ws.Range("F20").Value = 1227
$ws.Range("F21").Value = 554
$ws.Range("F22").Value = 1466
$ws.Range("F23").Value = 1330
$ws.Range("F24").Value = 1192
$ws.Range("F25").Value = 2268
$ws.Range("F26").Value = 355
$ws.Range("F28").Value = 1107
$ws.Range("F29").Value = 728
$ws.Range("F30").Value = 237
$ws.Range("F31").Value = 5299
$ws.Range("F32").Value = 284
$ws.Range("F33").Value = 1252
$ws.Range("F34").Value = 3705
$ws.Range("F36").Value = 1685
$ws.Range("F37").Value = 1074
$ws.Range("F38").Value = 158
$ws.Range("F39").Value = 960
$ws.Range("F40").Value = 390
$ws.Range("F41").Value = 1766
$ws.Range("F43").Value = 44
$ws.Range("F44").Value = 104
$ws.Range("F45").Value = 896
$ws.Range("F46").Value = 1046
$ws.Range("F47").Value = 510
$ws.Range("F48").Value = 193
$ws.Range("F49").Value = 193
$ws.Range("F51").Value = 77

$wb.Save()
